$d = $word.ActiveDocument

# --- Edit 1: "Class:" paragraph -> wrap existing text with spell/gram proofErr
#     markers and append a new "Investment" run, producing "Class:Investment"
#     (separate <w:r> from the trailing ":" run), matching the commit's
#     "GIC and investment" constructor rename.
$pClass = $d.Paragraphs(1)
$rClass = $pClass.Range

$classFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
  '<w:p w14:paraId="1AD42253" w14:textId="424C2EAD" w:rsidR="001D2D67" w:rsidRPr="001D2D67" w:rsidRDefault="001D2D67">' +
    '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidRPr="001D2D67"><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>C</w:t></w:r>' +
    '<w:r w:rsidRPr="001D2D67"><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>lass</w:t></w:r>' +
    '<w:r w:rsidRPr="001D2D67"><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>:</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Investment</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
  '</w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rClass.InsertXML($classFrag) | Out-Null

# --- Edit 2: empty "Responsibilities" paragraph -> drop the stray eastAsia
#     rFonts hint and give it a "Store info " run.
$pStore = $d.Paragraphs(4)
$rStore = $pStore.Range

$storeFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
  '<w:p w14:paraId="3D51B66F" w14:textId="0090EA92" w:rsidR="001D2D67" w:rsidRPr="001D2D67" w:rsidRDefault="001D2D67">' +
    '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Store info </w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rStore.InsertXML($storeFrag) | Out-Null

Write-Output ("Class paragraph: [" + $d.Paragraphs(1).Range.Text + "]")
Write-Output ("Store paragraph: [" + $d.Paragraphs(4).Range.Text + "]")
